$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "100K Resistor"
$ws.Range("B5").Value = 0.1
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "Stackpole Electronics, Inc"
$ws.Range("E5").Value = "RMCF0805JT100K"
$ws.Range("F5").Value = "RMCF0805JT100KCT-ND"
$ws.Range("G5").Value = "RES SMD 100K OHM 5% 1/8W 0805"

$ws.Hyperlinks.Add($ws.Range("I5"), "https://www.seielect.com/Catalog/SEI-RMCF_RMCP.pdf", "", "RMCF Series, Packaging Spec", "RMCF Series Datasheet")
$ws.Hyperlinks.Add($ws.Range("J5"), "https://www.seielect.com/catalog/SEI-Packaging.pdf", "", "", "Packaging Specs")
$ws.Hyperlinks.Add($ws.Range("H5"), "https://www.digikey.com/product-detail/en/stackpole-electronics-inc/RMCF0805JT100K/RMCF0805JT100KCT-ND/1942590", "", "", "Shopping")
